$wb = $excel.ActiveWorkbook

$wsSearch = $wb.Worksheets.Item("Search")
$wsSignIn = $wb.Worksheets.Item("SignIn")
$wsAcc    = $wb.Worksheets.Item("MyAcc_PInfo")

# ---------------------------------------------------------------------------
# Search sheet: add "SingleSearch" / "BrandSearch" columns (D & E)
# ---------------------------------------------------------------------------

# Header row (row 1) - reuse the existing header look from C1
$wsSearch.Range("D1").Value = "SingleSearch"
$wsSearch.Range("E1").Value = "BrandSearch"
$wsSearch.Range("C1").Copy()
$wsSearch.Range("D1:E1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Sample data row (row 2) - Honda brand search example
$wsSearch.Range("E2").Value = "Honda"
$wsSearch.Range("D2").Value = "IPI INDUSTRIES-HND1.2-000ID"

# D2 gets its own emphasis style: bold, dark grey Arial, wrapped & indented
$d2 = $wsSearch.Range("D2")
$d2.Font.Name = "Arial"
$d2.Font.Size = 9.9
$d2.Font.Bold = $true
$d2.Font.Color = 4473924
$d2.HorizontalAlignment = -4131   # xlLeft
$d2.VerticalAlignment = -4108     # xlCenter
$d2.WrapText = $true
$d2.IndentLevel = 1

# Row 2 grows a bit to fit the wrapped brand text
$wsSearch.Rows.Item(2).RowHeight = 26.4

# Column widths: new column D, plus a minor re-fit of A & C
$wsSearch.Columns.Item(1).ColumnWidth = 13.1640625
$wsSearch.Columns.Item(3).ColumnWidth = 11.83203125
$wsSearch.Columns.Item(4).ColumnWidth = 19.83203125

# ---------------------------------------------------------------------------
# Active tab / selection: "Search" becomes the active sheet (was "SignIn")
# ---------------------------------------------------------------------------
$wsSearch.Activate()
$wsSearch.Range("C12").Select()

# ---------------------------------------------------------------------------
# SignIn sheet: minor column re-fit (no content change)
# ---------------------------------------------------------------------------
$wsSignIn.Columns.Item(1).ColumnWidth = 28.83203125
$wsSignIn.Columns.Item(2).ColumnWidth = 10.83203125

# ---------------------------------------------------------------------------
# MyAcc_PInfo sheet: minor column re-fit (no content change)
# ---------------------------------------------------------------------------
$wsAcc.Columns.Item(3).ColumnWidth = 17.1640625
$wsAcc.Columns.Item(6).ColumnWidth = 15.83203125
$wsAcc.Columns.Item(9).ColumnWidth = 10.1640625
$wsAcc.Columns.Item(10).ColumnWidth = 12.5
$wsAcc.Columns.Item(12).ColumnWidth = 25.5
$wsAcc.Columns.Item(14).ColumnWidth = 15.83203125
